$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) format, used to reset
# style on cells we temporarily force to Text number-format so that
# decimal-looking strings (e.g. "150.00") are not coerced to numbers
# and lose their literal text representation (trailing zeros, etc).
$defaultStyleCell = $ws.Range("B2")

$ws.Range("D2").Value = "23.228.19"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.603.89"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = $defaultStyleCell.Style
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9988"
$ws.Range("D5").Style = $defaultStyleCell.Style
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.08"
$ws.Range("D6").Style = $defaultStyleCell.Style
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3778"
$ws.Range("D7").Style = $defaultStyleCell.Style
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.75"
$ws.Range("D8").Style = $defaultStyleCell.Style
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3634"
$ws.Range("D9").Style = $defaultStyleCell.Style
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.273"
$ws.Range("D10").Style = $defaultStyleCell.Style
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08127"
$ws.Range("D11").Style = $defaultStyleCell.Style
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9989"
$ws.Range("D12").Style = $defaultStyleCell.Style
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.92"
$ws.Range("D13").Style = $defaultStyleCell.Style
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.605"
$ws.Range("D14").Style = $defaultStyleCell.Style
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.426"
$ws.Range("D15").Style = $defaultStyleCell.Style
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001249"
$ws.Range("D16").Style = $defaultStyleCell.Style
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "1.603.66"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.91"
$ws.Range("D18").Style = $defaultStyleCell.Style
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06889"
$ws.Range("D19").Style = $defaultStyleCell.Style
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("D20").Style = $defaultStyleCell.Style
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.544"
$ws.Range("D21").Style = $defaultStyleCell.Style
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = $defaultStyleCell.Style
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.98"
$ws.Range("D23").Style = $defaultStyleCell.Style
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").Value = "23.235.78"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.008"
$ws.Range("D25").Style = $defaultStyleCell.Style
$ws.Range("E25").Value = "  +8.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.387"
$ws.Range("D26").Style = $defaultStyleCell.Style
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.25"
$ws.Range("D27").Style = $defaultStyleCell.Style
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.00"
$ws.Range("D28").Style = $defaultStyleCell.Style
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.20"
$ws.Range("D30").Style = $defaultStyleCell.Style
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.368"
$ws.Range("D31").Style = $defaultStyleCell.Style
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.775"
$ws.Range("D32").Style = $defaultStyleCell.Style
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "1.779.98"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9673"
$ws.Range("D34").Style = $defaultStyleCell.Style
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07513"
$ws.Range("D35").Style = $defaultStyleCell.Style
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02732"
$ws.Range("D36").Style = $defaultStyleCell.Style
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.26"
$ws.Range("D37").Style = $defaultStyleCell.Style
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2537"
$ws.Range("D38").Style = $defaultStyleCell.Style
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08808"
$ws.Range("D39").Style = $defaultStyleCell.Style
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.106"
$ws.Range("D40").Style = $defaultStyleCell.Style
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7128"
$ws.Range("D42").Style = $defaultStyleCell.Style
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.56"
$ws.Range("D43").Style = $defaultStyleCell.Style
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.66"
$ws.Range("D44").Style = $defaultStyleCell.Style
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6567"
$ws.Range("D45").Style = $defaultStyleCell.Style
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.321"
$ws.Range("D46").Style = $defaultStyleCell.Style
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.023"
$ws.Range("D47").Style = $defaultStyleCell.Style
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.36"
$ws.Range("D48").Style = $defaultStyleCell.Style
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07958"
$ws.Range("D49").Style = $defaultStyleCell.Style
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.206"
$ws.Range("D50").Style = $defaultStyleCell.Style
$ws.Range("E50").Value = "  -3.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.209"
$ws.Range("D51").Style = $defaultStyleCell.Style
$ws.Range("E51").Value = "  +0.83%  "
